# Daily update: decrement the "remaining days" counter (column E) for every
# data row. When a row's remaining count has hit 1 (i.e. the cycle is about
# to finish), reset it by adding back the "total days" value (column D) and
# advance the "start time" (column F, stored as yyyymmdd) by that same
# number of days - starting a fresh cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($row = 2; $row -le $lastRow; $row++) {
    $totalDays = $ws.Cells.Item($row, 4).Value2   # column D - total days
    $remaining = $ws.Cells.Item($row, 5).Value2   # column E - remaining days
    $startDate = $ws.Cells.Item($row, 6).Value2   # column F - start date (yyyymmdd)

    if ($remaining -eq $null -or $totalDays -eq $null) {
        continue
    }

    if ($remaining -le 1) {
        # Cycle finished: reset remaining and roll the start date forward.
        $newRemaining = $remaining - 1 + $totalDays

        $y = [Math]::Floor($startDate / 10000)
        $m = [Math]::Floor(($startDate % 10000) / 100)
        $d = $startDate % 100
        $dt = Get-Date -Year $y -Month $m -Day $d
        $newDt = $dt.AddDays($totalDays)
        $newStartDate = [int]($newDt.ToString("yyyyMMdd"))

        $ws.Cells.Item($row, 5).Value2 = $newRemaining
        $ws.Cells.Item($row, 6).Value2 = $newStartDate
    } else {
        $ws.Cells.Item($row, 5).Value2 = $remaining - 1
    }
}
